$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 5.5
$ws.Range("E2").Value = 0.09
$ws.Range("F2").Value = 331
$ws.Range("G2").Value = 0.0001
$ws.Range("I2").Value = 1000
$ws.Range("N2").Value = "[0.005]"
